$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.192.85'
$ws.Range("E2").Value = '  -1.37%  '

$ws.Range("D3").Value = '1.551.56'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.006'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.63'
$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3810'
$ws.Range("E7").Value = '  +3.98%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3255'
$ws.Range("E8").Value = '  -2.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.03'
$ws.Range("E9").Value = '  -9.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.125'
$ws.Range("E10").Value = '  -0.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07333'
$ws.Range("E11").Value = '  -1.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.23'
$ws.Range("E13").Value = '  -3.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.814'
$ws.Range("E14").Value = '  -3.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.751'
$ws.Range("E15").Value = '  -2.73%  '

$ws.Range("D16").Value = '1.563.61'
$ws.Range("E16").Value = '  -0.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001076'
$ws.Range("E17").Value = '  -3.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06736'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.50'
$ws.Range("E19").Value = '  -3.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.326'
$ws.Range("E21").Value = '  -1.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.12'
$ws.Range("E22").Value = '  -2.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.57'
$ws.Range("E23").Value = '  -4.78%  '

$ws.Range("D24").Value = '22.195.86'
$ws.Range("E24").Value = '  -1.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("E25").Value = '  -4.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.491'
$ws.Range("E26").Value = '  -5.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.64'
$ws.Range("E27").Value = '  -2.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.44'
$ws.Range("E28").Value = '  -1.23%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.942'
$ws.Range("E29").Value = '  -1.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.69'
$ws.Range("E30").Value = '  -1.35%  '

$ws.Range("D31").Value = '1.739.79'
$ws.Range("E31").Value = '  -0.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.041'
$ws.Range("E32").Value = '  -0.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.875'
$ws.Range("E33").Value = '  -5.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.855'
$ws.Range("E34").Value = '  -7.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.389'
$ws.Range("E35").Value = '  -5.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08260'
$ws.Range("E36").Value = '  -0.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02356'
$ws.Range("E37").Value = '  -3.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06205'
$ws.Range("E38").Value = '  -4.30%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.264'
$ws.Range("E39").Value = '  -3.06%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2154'
$ws.Range("E40").Value = '  -5.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.230'
$ws.Range("E41").Value = '  -4.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.96'
$ws.Range("E42").Value = '  -3.65%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.006'
$ws.Range("E43").Value = '  +0.42%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6011'
$ws.Range("E44").Value = '  -5.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.59'
$ws.Range("E45").Value = '  -2.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.740'
$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5853'
$ws.Range("E47").Value = '  -4.49%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.89'
$ws.Range("E48").Value = '  -1.84%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.971'
$ws.Range("E49").Value = '  -4.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.169'
$ws.Range("E50").Value = '  -4.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07046'
$ws.Range("E51").Value = '  -2.77%  '

